# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
# Swap the match-data (columns B:AD) between the two rows in each pair below,
# while leaving column A (the sequential row id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(21, 22),
    @(38, 39),
    @(156, 157),
    @(177, 178)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
